# Updates to Sections 3 & 4.
#
# - "Initialization Flag(s)" -> "Clean Flag(s)" text updates on the three
#   labelled boxes of the flag-register/flag-cache diagram.
# - The flag-cache box gains a second line ("Cache") on its own paragraph.
# - The "uninitialize" callout textbox is renamed to "clean" and is moved /
#   resized to sit over the narrower word.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Initialization Flags" -> "Clean Flags"
$flagsBox = $s.Shapes.Item("TextBox 77")
$flagsBox.TextFrame.TextRange.Text = "Clean Flags"

# "Initialization Flag Register File" -> "Clean Flag Register File"
$regFileBox = $s.Shapes.Item("Rectangle 79")
$regFileBox.TextFrame.TextRange.Text = "Clean Flag Register File"

# "Initialization Flag Cache" -> "Clean Flag" + new "Cache" paragraph
$cacheBox = $s.Shapes.Item("Rectangle 80")
$cacheBox.TextFrame.TextRange.Text = "Clean Flag" + [char]13 + "Cache"

# "uninitialize" -> "clean", with the textbox shrunk/repositioned to match
$cleanLabel = $s.Shapes.Item("TextBox 114")
$cleanLabel.TextFrame.TextRange.Text = "clean"
$cleanLabel.Left = 445.2451938503937
$cleanLabel.Top = 130.03685839370078
$cleanLabel.Width = 36.58
$cleanLabel.Height = 18.175748031496063
